$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.211.49"
$ws.Range("E2").Value = "  +6.33%  "
$ws.Range("D3").Value = "2.423.81"
$ws.Range("E3").Value = "  +2.56%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.29"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.46%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  +2.61%  "
$ws.Range("D9").Value = "2.457.52"
$ws.Range("E9").Value = "  +4.23%  "
$ws.Range("E10").Value = "  +6.21%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.355"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.15"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.96%  "
$ws.Range("E15").Value = "  +9.17%  "
$ws.Range("D16").Value = "2.959.73"
$ws.Range("E16").Value = "  +6.88%  "
$ws.Range("D17").Value = "62.653.09"
$ws.Range("E17").Value = "  +5.68%  "
$ws.Range("D18").Value = "2.437.91"
$ws.Range("E18").Value = "  +3.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.72"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +13.98%  "
$ws.Range("E24").Value = "  -0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.70"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "622.96"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +13.33%  "
$ws.Range("E27").Value = "  +10.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.42"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +8.34%  "
$ws.Range("D30").Value = "2.565.37"
$ws.Range("E31").Value = "  +3.28%  "
$ws.Range("E32").Value = "  +9.68%  "
$ws.Range("E33").Value = "  +6.50%  "
$ws.Range("E34").Value = "  +4.92%  "
$ws.Range("E35").Value = "  +6.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  +5.47%  "
$ws.Range("E38").Value = "  +2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.73"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.42"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.81%  "
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("E42").Value = "  +15.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.76"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.30%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "0.0₆0287"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.82"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("E48").Value = "  +3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.604"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.55%  "
$ws.Range("E51").Value = "  +4.21%  "
